$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 125000424
$ws.Range("I19").Value = 562.6667
$ws.Range("J19").Value = 500000000
$ws.Range("K19").Value = 562.6667
$ws.Range("L19").Value = 500000000
$ws.Range("M19").Value = -387.6667
$ws.Range("N19").Value = -500000350
$ws.Range("H93").Value = 46533.668
$ws.Range("J93").Value = 46533.668
$ws.Range("L93").Value = 46533.668
$ws.Range("N93").Value = -51525.668
$ws.Range("H116").Value = 2921.647
$ws.Range("I116").Value = 2709.5557
$ws.Range("J116").Value = 3160.25
$ws.Range("K116").Value = 2709.5557
$ws.Range("L116").Value = 3160.25
$ws.Range("M116").Value = 732.4443000000001
$ws.Range("N116").Value = -10044.25
$ws.Range("H129").Value = 1455.6041
$ws.Range("I129").Value = 629.4
$ws.Range("K129").Value = 1888.2
$ws.Range("M129").Value = 3111.8
$ws.Range("H132").Value = 4640.855
$ws.Range("I132").Value = 3580.5557
$ws.Range("J132").Value = 8457.933999999999
$ws.Range("K132").Value = 10741.6671
$ws.Range("L132").Value = 25373.802
$ws.Range("M132").Value = -8211.667099999999
$ws.Range("N132").Value = -30433.802

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1575
$ws.Range("J6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("N6").Value = -2346
$ws.Range("H26").Value = 471.66666
$ws.Range("I26").Value = 471.66666
$ws.Range("K26").Value = 471.66666
$ws.Range("M26").Value = -141.66666
$ws.Range("H32").Value = 9755.414000000001
$ws.Range("I32").Value = 7423.9326
$ws.Range("J32").Value = 30505.6
$ws.Range("K32").Value = 7423.9326
$ws.Range("L32").Value = 30505.6
$ws.Range("M32").Value = -7136.9326
$ws.Range("N32").Value = -31079.6
$ws.Range("H39").Value = 900
$ws.Range("I39").Value = 900
$ws.Range("K39").Value = 900
$ws.Range("M39").Value = -380
$ws.Range("H45").Value = 2617.077
$ws.Range("I45").Value = 1856
$ws.Range("J45").Value = 3834.8
$ws.Range("K45").Value = 1856
$ws.Range("L45").Value = 3834.8
$ws.Range("M45").Value = -1479
$ws.Range("N45").Value = -4588.8
$ws.Range("H61").Value = 1659.862
$ws.Range("I61").Value = 1527.1082
$ws.Range("K61").Value = 1527.1082
$ws.Range("M61").Value = -1315.1082
$ws.Range("H74").Value = 3097.8333
$ws.Range("I74").Value = 3857.1052
$ws.Range("J74").Value = 1786.3636
$ws.Range("K74").Value = 3857.1052
$ws.Range("L74").Value = 1786.3636
$ws.Range("M74").Value = -2983.1052
$ws.Range("N74").Value = -3534.3636
$ws.Range("H77").Value = 3097.8333
$ws.Range("I77").Value = 3857.1052
$ws.Range("J77").Value = 1786.3636
$ws.Range("K77").Value = 19285.526
$ws.Range("L77").Value = 8931.817999999999
$ws.Range("M77").Value = -14917.526
$ws.Range("N77").Value = -17667.818
$ws.Range("H98").Value = 38665.668
$ws.Range("J98").Value = 38665.668
$ws.Range("L98").Value = 38665.668
$ws.Range("N98").Value = -44655.668
$ws.Range("H132").Value = 3121.3494
$ws.Range("I132").Value = 1643.6964
$ws.Range("J132").Value = 6186.1113
$ws.Range("K132").Value = 4931.0892
$ws.Range("L132").Value = 18558.3339
$ws.Range("M132").Value = -2401.0892
$ws.Range("N132").Value = -23618.3339
$ws.Range("H136").Value = 1659.862
$ws.Range("I136").Value = 1527.1082
$ws.Range("K136").Value = 4581.3246
$ws.Range("M136").Value = -2031.3246

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2705.3076
$ws.Range("I99").Value = 2038.4286
$ws.Range("J99").Value = 3483.3333
$ws.Range("K99").Value = 2038.4286
$ws.Range("L99").Value = 3483.3333
$ws.Range("M99").Value = -540.4286
$ws.Range("N99").Value = -6479.3333
$ws.Range("H134").Value = 2439.4106
$ws.Range("I134").Value = 1488.5385
$ws.Range("J134").Value = 3263.5
$ws.Range("K134").Value = 4465.6155
$ws.Range("L134").Value = 9790.5
$ws.Range("M134").Value = -1930.6155
$ws.Range("N134").Value = -14860.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 532.9
$ws.Range("I19").Value = 532.9
$ws.Range("K19").Value = 532.9
$ws.Range("M19").Value = -362.9
$ws.Range("H24").Value = 532.9
$ws.Range("I24").Value = 532.9
$ws.Range("K24").Value = 532.9
$ws.Range("M24").Value = -362.9
$ws.Range("H44").Value = 3064
$ws.Range("I44").Value = 3064
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3064
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -2622
$ws.Range("N44").ClearContents()
$ws.Range("H93").Value = 8787.532999999999
$ws.Range("I93").Value = 4181.3
$ws.Range("J93").Value = 18000
$ws.Range("K93").Value = 4181.3
$ws.Range("L93").Value = 18000
$ws.Range("M93").Value = -2309.3
$ws.Range("N93").Value = -21744
$ws.Range("H99").Value = 3733.3333
$ws.Range("I99").Value = 3600
$ws.Range("K99").Value = 3600
$ws.Range("M99").Value = -2102
$ws.Range("H126").Value = 3733.3333
$ws.Range("I126").Value = 3600
$ws.Range("K126").Value = 10800
$ws.Range("M126").Value = -8330
$ws.Range("H132").Value = 2424.4424
$ws.Range("I132").Value = 1620
$ws.Range("K132").Value = 4860
$ws.Range("M132").Value = -2330

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 421.85715
$ws.Range("I2").Value = 30.5
$ws.Range("J2").Value = 715.375
$ws.Range("K2").Value = 183
$ws.Range("L2").Value = 4292.25
$ws.Range("M2").Value = -70
$ws.Range("N2").Value = -4518.25
$ws.Range("H38").Value = 218.75
$ws.Range("J38").Value = 116
$ws.Range("L38").Value = 348
$ws.Range("N38").Value = -1042
$ws.Range("H40").Value = 10163.5
$ws.Range("I40").Value = 125
$ws.Range("K40").Value = 500
$ws.Range("M40").Value = -431
$ws.Range("H69").Value = 1151
$ws.Range("J69").Value = 1857
$ws.Range("L69").Value = 5571
$ws.Range("N69").Value = -7193
$ws.Range("H72").Value = 1151
$ws.Range("J72").Value = 1857
$ws.Range("L72").Value = 16713
$ws.Range("N72").Value = -24825

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744
$ws.Range("H102").Value = 1994529
$ws.Range("I102").Value = 2756475.2
$ws.Range("J102").Value = 13468.8
$ws.Range("K102").Value = 2756475.2
$ws.Range("L102").Value = 13468.8
$ws.Range("M102").Value = -2754853.2
$ws.Range("N102").Value = -16712.8
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H130").Value = 31104.166
$ws.Range("J130").Value = 37849.152
$ws.Range("L130").Value = 37849.152
$ws.Range("N130").Value = -47889.152

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H68").Value = 2098.889
$ws.Range("I68").Value = 1707.6923
$ws.Range("K68").Value = 1707.6923
$ws.Range("M68").Value = -958.6922999999999
$ws.Range("H71").Value = 2098.889
$ws.Range("I71").Value = 1707.6923
$ws.Range("K71").Value = 8538.461499999999
$ws.Range("M71").Value = -4794.461499999999
$ws.Range("H112").Value = 22475
$ws.Range("J112").Value = 22475
$ws.Range("L112").Value = 22475
$ws.Range("N112").Value = -25429
$ws.Range("H136").Value = 1509.8788
$ws.Range("I136").Value = 1148
$ws.Range("J136").Value = 1894.375
$ws.Range("K136").Value = 3444
$ws.Range("L136").Value = 5683.125
$ws.Range("M136").Value = -894
$ws.Range("N136").Value = -10783.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H136").Value = 1987.9565
$ws.Range("I136").Value = 1521.6471
$ws.Range("J136").Value = 3309.1667
$ws.Range("K136").Value = 4564.9413
$ws.Range("L136").Value = 9927.500100000001
$ws.Range("M136").Value = -2014.9413
$ws.Range("N136").Value = -15027.5001
